# Mise a jour 1 du 25/04/2023
#
# After the paragraph ending in "... dans le vide." six new paragraphs are
# inserted: a relocated copy of the existing "Cardinaliter ..." paragraph,
# followed by five brand-new paragraphs ("Chaque  ligne ...", "Identifiant
# ...", "Id naurelle ...", "Id Artificielle inventer", "Id Composer = ").
# The trailing "_GoBack" bookmark moves from the end of the "... dans le
# vide." paragraph to the end of the new last paragraph ("Id Composer =
# "). The old "Cardinaliter ..." paragraph (which used to sit right after
# "... dans le vide.") is removed, together with one of the (now
# redundant) trailing empty paragraphs.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraph ending in "dans le vide." by searching its text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*dans le vide.*") {
        $target = $i
        break
    }
}

# Drop the "_GoBack" bookmark that currently sits at the end of that
# paragraph -- it will be recreated at the end of the new last paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Create six blank paragraphs right after the "dans le vide." paragraph.
$anchor = $d.Paragraphs.Item($target)
for ($i = 0; $i -lt 6; $i++) {
    $anchor.Range.InsertParagraphAfter() | Out-Null
}

$p = $target + 1

# 1) Relocated copy of the "Cardinaliter ..." paragraph.
$xml = "<w:p $wNs>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Cardinaliter</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> exemple </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/><w:r><w:t>de acte</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> national Père, </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Mere</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> unique mais enfant plusieurs.</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$p++

# 2) "Chaque  ligne represente une occurrence."
$xml = "<w:p $wNs>" +
    "<w:proofErr w:type=`"gramStart`"/><w:r><w:t>Chaque  ligne</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>represente</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> une occurrence.</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$p++

# 3) "Identifiant a partir de xa on peutuj avoir toute les renseignement sur vous."
$xml = "<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">Identifiant </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> partir de </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>xa</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> on </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>peutuj</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> avoir toute les renseignement sur vous.</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$p++

# 4) "Id naurelle = il y'a xa dans les propriete"
$xml = "<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">Id </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>naurelle</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> = il y’a </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>xa</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> dans les </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>propriete</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$p++

# 5) "Id Artificielle inventer"
$xml = "<w:p $wNs><w:r><w:t>Id Artificielle inventer</w:t></w:r></w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$p++

# 6) "Id Composer = " + the relocated _GoBack bookmark.
$xml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Id Composer = </w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$d.Paragraphs.Item($p).Range.InsertXML($xml) | Out-Null
$lastNew = $p

# Remove the old "Cardinaliter ..." paragraph that used to follow the
# "dans le vide." paragraph directly (it now lives earlier in the
# document, right after "dans le vide.").
$cardIndex = $null
for ($i = $lastNew + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Cardinaliter*") {
        $cardIndex = $i
        break
    }
}
if ($cardIndex) {
    $d.Paragraphs.Item($cardIndex).Range.Delete() | Out-Null
}

# Remove one of the (now one-too-many) trailing empty paragraphs. Walk
# forward and stop before the very last paragraph -- Word (and this COM
# shim) will not delete the document's final paragraph mark.
for ($i = $lastNew + 1; $i -lt $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text -eq "`r") {
        $par.Range.Delete() | Out-Null
        break
    }
}
